$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" updates ---
$ws1 = $wb.Worksheets.Item("Schedule")

$ws1.Range("E3").Value = 441.157782
$ws1.Range("F3").Value = 29.17710198412698

$ws1.Range("A4").Value = 46039.33333333334
$ws1.Range("B4").Value = 46039.52083333334
$ws1.Range("C4").Value = 4.5
$ws1.Range("D4").Value = 17.01
$ws1.Range("E4").Value = 132.51114825
$ws1.Range("F4").Value = 7.790190961199297

# New row 5
$ws1.Range("A5").Value = 46039.60416666666
$ws1.Range("B5").Value = 46039.91666666666
$ws1.Range("C5").Value = 7.5
$ws1.Range("D5").Value = 28.35
$ws1.Range("E5").Value = 91.56559425
$ws1.Range("F5").Value = 3.229826957671958
$ws1.Range("A5").NumberFormat = $ws1.Range("A4").NumberFormat
$ws1.Range("B5").NumberFormat = $ws1.Range("B4").NumberFormat

# --- Sheet "Detailed" updates ---
$ws2 = $wb.Worksheets.Item("Detailed")

$ws2.Range("B41").Value = 22.01959
$ws2.Range("B42").Value = 22.01959
$ws2.Range("B43").Value = 69.68688
$ws2.Range("C43").Value = "historical"
$ws2.Range("B44").Value = 62.33685
$ws2.Range("C44").Value = "historical"
$ws2.Range("B45").Value = 62.33685
$ws2.Range("B48").Value = 57.09
$ws2.Range("B49").Value = 57.06003
$ws2.Range("B51").Value = 56.98
$ws2.Range("B52").Value = 56.98
$ws2.Range("B53").Value = 56.98
$ws2.Range("B59").Value = 58.01829
$ws2.Range("B60").Value = 56.98
$ws2.Range("B61").Value = 57.06017
$ws2.Range("B64").Value = 36.06029
$ws2.Range("B65").Value = 24.52555
$ws2.Range("E65").Value = "OFF"
$ws2.Range("B66").Value = 17.45481
$ws2.Range("B67").Value = 35.88
$ws2.Range("B68").Value = 27.45385
$ws2.Range("B69").Value = 19.65591
$ws2.Range("B70").Value = 0.7
$ws2.Range("B71").Value = 0.7
$ws2.Range("B72").Value = 0.51003
$ws2.Range("B73").Value = 11.48427
$ws2.Range("B74").Value = 22.07
$ws2.Range("B75").Value = 34.01
$ws2.Range("E75").Value = "OFF"
$ws2.Range("B76").Value = 34.01
$ws2.Range("E76").Value = "OFF"
$ws2.Range("B77").Value = 22.95651
$ws2.Range("E77").Value = "OFF"
$ws2.Range("B78").Value = 22.07
$ws2.Range("E78").Value = "OFF"
$ws2.Range("B79").Value = 7.72339
$ws2.Range("B80").Value = -2.21718
$ws2.Range("B81").Value = 8.834669999999999
$ws2.Range("B82").Value = 30.36368
$ws2.Range("B83").Value = 0.3657
$ws2.Range("B84").Value = -6
$ws2.Range("B85").Value = -3.03666
$ws2.Range("B86").Value = -6.20523
$ws2.Range("B87").Value = -6
$ws2.Range("B88").Value = 0.0109
$ws2.Range("B89").Value = 10.0727
$ws2.Range("E89").Value = "ON"
$ws2.Range("B90").Value = 13.59537
$ws2.Range("E90").Value = "ON"
$ws2.Range("B91").Value = 29.64754
$ws2.Range("E91").Value = "ON"
$ws2.Range("B92").Value = 8.43773
$ws2.Range("E92").Value = "ON"
$ws2.Range("B93").Value = 8.320819999999999
$ws2.Range("E93").Value = "ON"
$ws2.Range("B94").Value = 30.04051
$ws2.Range("B96").Value = 57.06002
